# Daily attendance processing - 2026-01-08 11:58:41
# Reverse the order of the comma-separated "Recorded By" entries (column G)
# for every data row that has more than one recorder listed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $value = $cell.Text

    if ($value.Contains(",")) {
        $parts = $value -split ",\s*"
        $n = $parts.Length
        $rev = @()
        for ($i = $n - 1; $i -ge 0; $i--) {
            $rev += $parts[$i]
        }
        $joined = [string]::Join(", ", $rev)
        $cell.Value = $joined
    }
}
